# Applies the "automatic update" diff against rows 2-12 of the Artfynd sheet.
#
# The 11 species-observation records that occupy sheet rows 2-12 get
# reshuffled into a new row order, and the "Taxonsorteringsordning"
# (column B) value for each record is refreshed to a new id. Columns
# C, I, K, P, S, T, U, V, W, Y, AA, AD, AE, AG, AT, AW, AX, AY are identical
# across every one of these rows (same observation batch/date/observer), so
# they don't need to be touched - only the species-identifying columns
# (A, B, D, E, F, G, H, Q, R) and the one-off "Aktivitet" note in column M
# (which travels with the "Tretåig hackspett" record) are written.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New sheet row -> the record (keyed by its old row) that now lives there,
# plus that record's refreshed column-B id.
$rowMap = @(
    @{ NewRow = 2;  OldRow = 5;  B = 78740 },
    @{ NewRow = 3;  OldRow = 12; B = 56430 },
    @{ NewRow = 4;  OldRow = 4;  B = 89571 },
    @{ NewRow = 5;  OldRow = 6;  B = 85850 },
    @{ NewRow = 6;  OldRow = 8;  B = 89993 },
    @{ NewRow = 7;  OldRow = 10; B = 78739 },
    @{ NewRow = 8;  OldRow = 7;  B = 89553 },
    @{ NewRow = 9;  OldRow = 2;  B = 89517 },
    @{ NewRow = 10; OldRow = 11; B = 89834 },
    @{ NewRow = 11; OldRow = 9;  B = 78713 },
    @{ NewRow = 12; OldRow = 3;  B = 89517 }
)

$cols = @("A", "D", "E", "F", "G", "H", "Q", "R")

# Snapshot the "before" values for the columns we touch, for every row,
# before any writes happen (writes below would otherwise clobber a value
# that a later row still needs to read as its source).
$before = @{}
foreach ($entry in $rowMap) {
    $oldRow = $entry.OldRow
    if (-not $before.ContainsKey($oldRow)) {
        $rowVals = @{}
        foreach ($col in $cols) {
            $rowVals[$col] = $ws.Range("$col$oldRow").Value2
        }
        $rowVals["M"] = $ws.Range("M$oldRow").Value2
        $before[$oldRow] = $rowVals
    }
}

foreach ($entry in $rowMap) {
    $newRow = $entry.NewRow
    $src = $before[$entry.OldRow]

    foreach ($col in $cols) {
        $ws.Range("$col$newRow").Value2 = $src[$col]
    }
    $ws.Range("B$newRow").Value2 = $entry.B

    if ($null -eq $src["M"]) {
        $ws.Range("M$newRow").ClearContents()
    } else {
        $ws.Range("M$newRow").Value2 = $src["M"]
    }
}
